$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 keeps its existing "plain number" style (cellXfs index 1, the same
# unused style already sitting on rows 5-8/13-20). Grab that style via a
# format-only paste instead of re-declaring font/fill properties, so the
# engine reuses the existing xf entry rather than minting a new one.
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4122) | Out-Null
$ws.Range("A3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# New barcode values (Trendyol / Morhipo)
$ws.Range("A2").Value = 8699490222850
$ws.Range("A3").Value = 8699490221419

# Row 2 already carries ht=16.5; give row 3 the same visual height.
$ws.Rows.Item(3).RowHeight = 16.5

# Drop the old leftover rows (4-20) that held no-longer-needed placeholder cells.
$ws.Rows.Item(4).Resize(17).Delete() | Out-Null

# Widen column A to fit the longer barcode values.
$ws.Columns.Item(1).ColumnWidth = 32

# Move the active selection.
$ws.Range("F5").Select() | Out-Null

# Paper size for the single print area.
$ws.PageSetup.PaperSize = 9
